# Updates the crypto price table to reflect a refreshed data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the cells we touch so Excel does not auto-convert
# numeric-looking strings (prices) or percentage-looking strings (volumes)
# into real numbers - the source data keeps everything as literal text.
# (Kept as two separate Range calls - a single comma-joined multi-area range
# does not reliably propagate NumberFormat across all of its areas.)
$bcRange = $ws.Range("B6:C18")
$deRange = $ws.Range("D2:E49")
$bcRange.NumberFormat = "@"
$deRange.NumberFormat = "@"

# Rows 6-18: the coin list shifted up by one position (GateToken, which was
# last, rotated to the top of this block) and prices/volumes were refreshed.
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "4.298"
$ws.Range("E6").Value = "1.28%"

$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "1.627"
$ws.Range("E7").Value = "1.71%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9149"
$ws.Range("E8").Value = "1.24%"

$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "2.444"
$ws.Range("E9").Value = "0.74%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1212"
$ws.Range("E10").Value = "9.48%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1824"
$ws.Range("E11").Value = "2.40%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09133"
$ws.Range("E12").Value = "0.20%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.04216"
$ws.Range("E13").Value = "0.25%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.1051"
$ws.Range("E14").Value = "-0.22%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001263"
$ws.Range("E15").Value = "0.48%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005871"
$ws.Range("E16").Value = "3.80%"

$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "0.007509"
$ws.Range("E17").Value = "1,903.11%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.342"
$ws.Range("E18").Value = "-0.18%"

# Remaining price / volume(1h) refreshes (no coin reordering)
$ws.Range("D2").Value = "308.37"
$ws.Range("E2").Value = "-0.03%"

$ws.Range("E3").Value = "0.48%"

$ws.Range("E4").Value = "1.99%"

$ws.Range("D5").Value = "0.07683"
$ws.Range("E5").Value = "0.61%"

$ws.Range("D20").Value = "7.405"
$ws.Range("E20").Value = "13.19%"

$ws.Range("D21").Value = "0.1381"
$ws.Range("E21").Value = "1.24%"

$ws.Range("D22").Value = "0.2711"
$ws.Range("E22").Value = "-4.18%"

$ws.Range("D23").Value = "0.04022"
$ws.Range("E23").Value = "-1.05%"

$ws.Range("D24").Value = "0.001261"
$ws.Range("E24").Value = "2.61%"

$ws.Range("D25").Value = "0.004271"
$ws.Range("E25").Value = "3.43%"

$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").Value = "-0.02%"

$ws.Range("D38").Value = "0.02502"
$ws.Range("E38").Value = "3.88%"

$ws.Range("D39").Value = "0.05316"
$ws.Range("E39").Value = "2.57%"

$ws.Range("D40").Value = "0.007852"
$ws.Range("E40").Value = "0.68%"

$ws.Range("D41").Value = "0.1315"
$ws.Range("E41").Value = "0.80%"

$ws.Range("D42").Value = "0.006503"
$ws.Range("E42").Value = "-7.73%"

$ws.Range("D43").Value = "0.001862"
$ws.Range("E43").Value = "-4.64%"

$ws.Range("D44").Value = "0.008032"
$ws.Range("E44").Value = "-8.89%"

$ws.Range("D45").Value = "0.3041"
$ws.Range("E45").Value = "-8.79%"

$ws.Range("D46").Value = "0.00006712"
$ws.Range("E46").Value = "-3.50%"

$ws.Range("D48").Value = "0.2795"
$ws.Range("E48").Value = "808.99%"

$ws.Range("D49").Value = "0.003101"
$ws.Range("E49").Value = "-26.18%"

# Restore the default (un-styled) cell format now that the text values are set,
# matching the workbook's original styling for these cells.
$bcRange.Style = "Normal"
$deRange.Style = "Normal"
